# Saldo_guide.xlsx update: refresh extraction date (2024-05-23 -> 2024-05-24)
# across the whole sheet, rename the sheet to match the new timestamp, and
# correct the Saldo Previsto / Vl. Total figures for row 109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new extraction run.
$ws.Name = "IClientBalance-20240524-092026-"

# Column G ("Dt. Referencia") holds the reference date serial for every data
# row (2..257). Bump it from 45435 (2024-05-23) to 45436 (2024-05-24).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
For ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45436
}

# Row 109 ("Saldo Previsto" / "Vl. Total") was corrected to 12221.02.
$ws.Cells.Item(109, 4).Value = 12221.02
$ws.Cells.Item(109, 8).Value = 12221.02
